$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047872834555001
$ws.Range("D2").Value = 1.046888489374791
$ws.Range("E2").Value = 1.051543531943921
$ws.Range("F2").Value = 1.046749836017703
$ws.Range("I2").Value = 1.042026306375352
$ws.Range("J2").Value = 1.052919589952319
$ws.Range("K2").Value = 1.049652664333931
$ws.Range("L2").Value = 1.054294741769581
$ws.Range("M2").Value = 1.04951439906931
$ws.Range("N2").Value = 1.054414856541616
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049344227391863
$ws.Range("D3").Value = 1.047995766124216
$ws.Range("E3").Value = 1.052965283923206
$ws.Range("F3").Value = 1.048847537169125
$ws.Range("I3").Value = 1.042488008664078
$ws.Range("J3").Value = 1.054037069076656
$ws.Range("K3").Value = 1.0505709073597
$ws.Range("L3").Value = 1.055527595626825
$ws.Range("M3").Value = 1.051420470117733
$ws.Range("N3").Value = 1.055533922614486
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.05029399537022
$ws.Range("D4").Value = 1.048710217396566
$ws.Range("E4").Value = 1.053883294080731
$ws.Range("F4").Value = 1.050202220319866
$ws.Range("I4").Value = 1.04278437326281
$ws.Range("J4").Value = 1.054757430456878
$ws.Range("K4").Value = 1.051162458568309
$ws.Range("L4").Value = 1.056322859395637
$ws.Range("M4").Value = 1.052650791776003
$ws.Range("N4").Value = 1.056255306990495
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.050692732305291
$ws.Range("D5").Value = 1.049010093383707
$ws.Range("E5").Value = 1.054268765998191
$ws.Range("F5").Value = 1.050771109384476
$ws.Range("I5").Value = 1.042908397068182
$ws.Range("J5").Value = 1.055059626824334
$ws.Range("K5").Value = 1.051410527589777
$ws.Range("L5").Value = 1.056656604366423
$ws.Range("M5").Value = 1.053167310434213
$ws.Range("N5").Value = 1.056557932511455
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.050759650298206
$ws.Range("D6").Value = 1.049060415971461
$ws.Range("E6").Value = 1.054333461749694
$ws.Range("F6").Value = 1.050866592581346
$ws.Range("I6").Value = 1.042929188018482
$ws.Range("J6").Value = 1.055110329421878
$ws.Range("K6").Value = 1.051452143354615
$ws.Range("L6").Value = 1.056712607623688
$ws.Range("M6").Value = 1.053253995151153
$ws.Range("N6").Value = 1.056608707112503
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050299325443566
$ws.Range("D7").Value = 1.048714226229198
$ws.Range("E7").Value = 1.053888446569648
$ws.Range("F7").Value = 1.05020982425467
$ws.Range("I7").Value = 1.042786032701737
$ws.Range("J7").Value = 1.054761470938235
$ws.Range("K7").Value = 1.051165775703345
$ws.Range("L7").Value = 1.05632732119607
$ws.Range("M7").Value = 1.052657696286028
$ws.Range("N7").Value = 1.056259353209799
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048370585555703
$ws.Range("D8").Value = 1.047263123378116
$ws.Range("E8").Value = 1.052024430813081
$ws.Range("F8").Value = 1.047459328467891
$ws.Range("I8").Value = 1.042182837465354
$ws.Range("J8").Value = 1.053297815874781
$ws.Range("K8").Value = 1.049963533954241
$ws.Range("L8").Value = 1.054711908208445
$ws.Range("M8").Value = 1.050159203575695
$ws.Range("N8").Value = 1.054793619588271
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044953656230705
$ws.Range("D9").Value = 1.044690231555894
$ws.Range("E9").Value = 1.048724378137043
$ws.Range("F9").Value = 1.042591281839346
$ws.Range("I9").Value = 1.041101475051429
$ws.Range("J9").Value = 1.050697450035466
$ws.Range("K9").Value = 1.047824718867227
$ws.Range("L9").Value = 1.05184599851778
$ws.Range("M9").Value = 1.045732506870338
$ws.Range("N9").Value = 1.052189560931198
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042662787026061
$ws.Range("D10").Value = 1.042963868471546
$ws.Range("E10").Value = 1.046513386140269
$ws.Range("F10").Value = 1.039330355189932
$ws.Range("I10").Value = 1.040367920446134
$ws.Range("J10").Value = 1.048949087565126
$ws.Range("K10").Value = 1.046384769095094
$ws.Range("L10").Value = 1.049921837428482
$ws.Range("M10").Value = 1.042764096845558
$ws.Range("N10").Value = 1.050438715585611
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041667605970281
$ws.Range("D11").Value = 1.042213602826281
$ws.Range("E11").Value = 1.045553270964612
$ws.Range("F11").Value = 1.03791437269429
$ws.Range("I11").Value = 1.040047230088561
$ws.Range("J11").Value = 1.048188403951734
$ws.Range("K11").Value = 1.045757823341036
$ws.Range("L11").Value = 1.049085316539662
$ws.Range("M11").Value = 1.041474393829656
$ws.Range("N11").Value = 1.049676951714236
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041297455044862
$ws.Range("D12").Value = 1.041934500686957
$ws.Range("E12").Value = 1.045196218115678
$ws.Range("F12").Value = 1.037387792144943
$ws.Range("I12").Value = 1.039927647533707
$ws.Range("J12").Value = 1.047905296509828
$ws.Range("K12").Value = 1.045524423080515
$ws.Range("L12").Value = 1.048774081974724
$ws.Range("M12").Value = 1.040994663723439
$ws.Range("N12").Value = 1.04939344222729
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041376876341529
$ws.Range("D13").Value = 1.041994388191032
$ws.Range("E13").Value = 1.045272826506293
$ws.Range("F13").Value = 1.037500773972424
$ws.Range("I13").Value = 1.039953319464986
$ws.Range("J13").Value = 1.04796604932683
$ws.Range("K13").Value = 1.045574512104347
$ws.Range("L13").Value = 1.048840866253599
$ws.Range("M13").Value = 1.041097598444881
$ws.Range("N13").Value = 1.04945428132026
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041637019378901
$ws.Range("D14").Value = 1.042190540775644
$ws.Range("E14").Value = 1.045523765567463
$ws.Range("F14").Value = 1.037870858200377
$ws.Range("I14").Value = 1.040037354847567
$ws.Range("J14").Value = 1.04816501359703
$ws.Range("K14").Value = 1.045738541153878
$ws.Range("L14").Value = 1.049059600330186
$ws.Range("M14").Value = 1.041434753125827
$ws.Range("N14").Value = 1.049653528142545
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041797235934769
$ws.Range("D15").Value = 1.042311340935188
$ws.Range("E15").Value = 1.045678320964009
$ws.Range("F15").Value = 1.038098796225549
$ws.Range("I15").Value = 1.040089070210749
$ws.Range("J15").Value = 1.048287528112354
$ws.Range("K15").Value = 1.04583953509196
$ws.Range("L15").Value = 1.049194301253321
$ws.Range("M15").Value = 1.041642394893418
$ws.Range("N15").Value = 1.049776216642534
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042728765004484
$ws.Range("D16").Value = 1.043013602694454
$ws.Range("E16").Value = 1.046577047091976
$ws.Range("F16").Value = 1.039424243420023
$ws.Range("I16").Value = 1.040389138815502
$ws.Range("J16").Value = 1.048999494264515
$ws.Range("K16").Value = 1.046426304289906
$ws.Range("L16").Value = 1.049977283149006
$ws.Range("M16").Value = 1.042849596444216
$ws.Range("N16").Value = 1.050489193868294
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043312217709085
$ws.Range("D17").Value = 1.043453373610083
$ws.Range("E17").Value = 1.047140052304061
$ws.Range("F17").Value = 1.040254579296006
$ws.Range("I17").Value = 1.040576542350399
$ws.Range("J17").Value = 1.049445112541534
$ws.Range("K17").Value = 1.046793442915415
$ws.Range("L17").Value = 1.050467523555612
$ws.Range("M17").Value = 1.043605659069221
$ws.Range("N17").Value = 1.05093544497438
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043652226021259
$ws.Range("D18").Value = 1.043709620851523
$ws.Range("E18").Value = 1.047468180031845
$ws.Range("F18").Value = 1.040738517790792
$ws.Range("I18").Value = 1.04068555708599
$ws.Range("J18").Value = 1.049704684483505
$ws.Range("K18").Value = 1.047007257504163
$ws.Range("L18").Value = 1.050753150336644
$ws.Range("M18").Value = 1.044046238222928
$ws.Range("N18").Value = 1.051195385538282
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04376810790618
$ws.Range("D19").Value = 1.043796950069241
$ws.Range("E19").Value = 1.047580018727699
$ws.Range("F19").Value = 1.040903464167184
$ws.Range("I19").Value = 1.040722678469507
$ws.Range("J19").Value = 1.049793132804408
$ws.Range("K19").Value = 1.047080106828638
$ws.Range("L19").Value = 1.050850487366455
$ws.Range("M19").Value = 1.044196393845644
$ws.Range("N19").Value = 1.051283959465946
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043249650889048
$ws.Range("D20").Value = 1.04340621768971
$ws.Range("E20").Value = 1.047079674530605
$ws.Range("F20").Value = 1.040165531836351
$ws.Range("I20").Value = 1.040556466227925
$ws.Range("J20").Value = 1.049397338157919
$ws.Range("K20").Value = 1.046754086727941
$ws.Range("L20").Value = 1.050414958798739
$ws.Range("M20").Value = 1.04352458423997
$ws.Range("N20").Value = 1.050887602745661
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041560427542965
$ws.Range("D21").Value = 1.0421327903938
$ws.Range("E21").Value = 1.045449882047333
$ws.Range("F21").Value = 1.037761894970274
$ws.Range("I21").Value = 1.040012621365158
$ws.Range("J21").Value = 1.048106438999982
$ws.Range("K21").Value = 1.045690253237292
$ws.Range("L21").Value = 1.048995202856962
$ws.Range("M21").Value = 1.041335488338817
$ws.Range("N21").Value = 1.049594870362851
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040495467556831
$ws.Range("D22").Value = 1.041329701381949
$ws.Range("E22").Value = 1.044422711811401
$ws.Range("F22").Value = 1.036247024732196
$ws.Range("I22").Value = 1.039667998300064
$ws.Range("J22").Value = 1.047291579371914
$ws.Range("K22").Value = 1.045018339071452
$ws.Range("L22").Value = 1.04809956983512
$ws.Range("M22").Value = 1.039955189144558
$ws.Range("N22").Value = 1.048778853540645
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04106030008183
$ws.Range("D23").Value = 1.041755667731154
$ws.Range("E23").Value = 1.044967470771608
$ws.Range("F23").Value = 1.037050436116802
$ws.Range("I23").Value = 1.039850945709842
$ws.Range("J23").Value = 1.047723860521692
$ws.Range("K23").Value = 1.045374824376816
$ws.Range("L23").Value = 1.048574647548831
$ws.Range("M23").Value = 1.040687291460085
$ws.Range("N23").Value = 1.049211748579239
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043277923099047
$ws.Range("D24").Value = 1.043427526235193
$ws.Range("E24").Value = 1.047106957462302
$ws.Range("F24").Value = 1.040205769731395
$ws.Range("I24").Value = 1.040565538673738
$ws.Range("J24").Value = 1.049418926425866
$ws.Range("K24").Value = 1.046771871118277
$ws.Range("L24").Value = 1.050438711545185
$ws.Range("M24").Value = 1.043561219760933
$ws.Range("N24").Value = 1.050909221671425
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045839242849683
$ws.Range("D25").Value = 1.045357310251699
$ws.Range("E25").Value = 1.049579408432722
$ws.Range("F25").Value = 1.043852440815402
$ws.Range("I25").Value = 1.041383244544835
$ws.Range("J25").Value = 1.051372274541286
$ws.Range("K25").Value = 1.048380103470442
$ws.Range("L25").Value = 1.052589253451003
$ws.Range("M25").Value = 1.046879876760653
$ws.Range("N25").Value = 1.05286534376522
